$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja4")

# Rows 57-59 / 60-62 (financed-interest block): TNA bancaria (H) increased;
# J (interes del banco), K (precio final financiado) and L (TNA total) are
# plain cached values (no formulas in the sheet) so they are updated to match.
$ws.Range("H57").Value = 0.65
$ws.Range("J57").Value = 32.054794520547944
$ws.Range("K57").Value = 146.44876712328767
$ws.Range("L57").Value = 0.9418777777777777
$ws.Range("H58").Value = 0.71
$ws.Range("J58").Value = 52.520547945205479
$ws.Range("K58").Value = 169.14528767123286
$ws.Range("L58").Value = 0.93474185185185177
$ws.Range("H59").Value = 0.74
$ws.Range("J59").Value = 72.986301369863014
$ws.Range("K59").Value = 191.84180821917809
$ws.Range("L59").Value = 0.93117388888888908
$ws.Range("H60").Value = 0.65
$ws.Range("J60").Value = 32.054794520547944
$ws.Range("K60").Value = 146.44876712328767
$ws.Range("L60").Value = 0.9418777777777777
$ws.Range("H61").Value = 0.71
$ws.Range("J61").Value = 52.520547945205479
$ws.Range("K61").Value = 169.14528767123286
$ws.Range("L61").Value = 0.93474185185185177
$ws.Range("H62").Value = 0.74
$ws.Range("J62").Value = 72.986301369863014
$ws.Range("K62").Value = 191.84180821917809
$ws.Range("L62").Value = 0.93117388888888908

# Rows 76-110 (Tasa 0% block): Costo de Procesamiento (I) increased;
# K and L recalculated accordingly (H=0 / J=0 unchanged).
$ws.Range("I76").Value = 12
$ws.Range("K76").Value = 112.00000000000001
$ws.Range("L76").Value = 0.73000000000000065
$ws.Range("I77").Value = 16
$ws.Range("K77").Value = 115.99999999999999
$ws.Range("L77").Value = 0.64888888888888863
$ws.Range("I78").Value = 21
$ws.Range("K78").Value = 121
$ws.Range("L78").Value = 0.63874999999999982
$ws.Range("I79").Value = 26
$ws.Range("K79").Value = 126
$ws.Range("L79").Value = 0.63266666666666671
$ws.Range("I80").Value = 31
$ws.Range("K80").Value = 131
$ws.Range("L80").Value = 0.62861111111111123
$ws.Range("I81").Value = 12
$ws.Range("K81").Value = 112.00000000000001
$ws.Range("L81").Value = 0.73000000000000065
$ws.Range("I82").Value = 16
$ws.Range("K82").Value = 115.99999999999999
$ws.Range("L82").Value = 0.64888888888888863
$ws.Range("I83").Value = 21
$ws.Range("K83").Value = 121
$ws.Range("L83").Value = 0.63874999999999982
$ws.Range("I84").Value = 26
$ws.Range("K84").Value = 126
$ws.Range("L84").Value = 0.63266666666666671
$ws.Range("I85").Value = 31
$ws.Range("K85").Value = 131
$ws.Range("L85").Value = 0.62861111111111123
$ws.Range("I86").Value = 12
$ws.Range("K86").Value = 112.00000000000001
$ws.Range("L86").Value = 0.73000000000000065
$ws.Range("I87").Value = 16
$ws.Range("K87").Value = 115.99999999999999
$ws.Range("L87").Value = 0.64888888888888863
$ws.Range("I88").Value = 21
$ws.Range("K88").Value = 121
$ws.Range("L88").Value = 0.63874999999999982
$ws.Range("I89").Value = 26
$ws.Range("K89").Value = 126
$ws.Range("L89").Value = 0.63266666666666671
$ws.Range("I90").Value = 31
$ws.Range("K90").Value = 131
$ws.Range("L90").Value = 0.62861111111111123
$ws.Range("I91").Value = 12
$ws.Range("K91").Value = 112.00000000000001
$ws.Range("L91").Value = 0.73000000000000065
$ws.Range("I92").Value = 16
$ws.Range("K92").Value = 115.99999999999999
$ws.Range("L92").Value = 0.64888888888888863
$ws.Range("I93").Value = 21
$ws.Range("K93").Value = 121
$ws.Range("L93").Value = 0.63874999999999982
$ws.Range("I94").Value = 26
$ws.Range("K94").Value = 126
$ws.Range("L94").Value = 0.63266666666666671
$ws.Range("I95").Value = 31
$ws.Range("K95").Value = 131
$ws.Range("L95").Value = 0.62861111111111123
$ws.Range("I96").Value = 12
$ws.Range("K96").Value = 112.00000000000001
$ws.Range("L96").Value = 0.73000000000000065
$ws.Range("I97").Value = 16
$ws.Range("K97").Value = 115.99999999999999
$ws.Range("L97").Value = 0.64888888888888863
$ws.Range("I98").Value = 21
$ws.Range("K98").Value = 121
$ws.Range("L98").Value = 0.63874999999999982
$ws.Range("I99").Value = 26
$ws.Range("K99").Value = 126
$ws.Range("L99").Value = 0.63266666666666671
$ws.Range("I100").Value = 31
$ws.Range("K100").Value = 131
$ws.Range("L100").Value = 0.62861111111111123
$ws.Range("I101").Value = 12
$ws.Range("K101").Value = 112.00000000000001
$ws.Range("L101").Value = 0.73000000000000065
$ws.Range("I102").Value = 16
$ws.Range("K102").Value = 115.99999999999999
$ws.Range("L102").Value = 0.64888888888888863
$ws.Range("I103").Value = 21
$ws.Range("K103").Value = 121
$ws.Range("L103").Value = 0.63874999999999982
$ws.Range("I104").Value = 26
$ws.Range("K104").Value = 126
$ws.Range("L104").Value = 0.63266666666666671
$ws.Range("I105").Value = 31
$ws.Range("K105").Value = 131
$ws.Range("L105").Value = 0.62861111111111123
$ws.Range("I106").Value = 12
$ws.Range("K106").Value = 112.00000000000001
$ws.Range("L106").Value = 0.73000000000000065
$ws.Range("I107").Value = 16
$ws.Range("K107").Value = 115.99999999999999
$ws.Range("L107").Value = 0.64888888888888863
$ws.Range("I108").Value = 21
$ws.Range("K108").Value = 121
$ws.Range("L108").Value = 0.63874999999999982
$ws.Range("I109").Value = 26
$ws.Range("K109").Value = 126
$ws.Range("L109").Value = 0.63266666666666671
$ws.Range("I110").Value = 31
$ws.Range("K110").Value = 131
$ws.Range("L110").Value = 0.62861111111111123

# Rows 121-125 (Mixto block): TNA bancaria (H) increased; J, K, L recalculated.
$ws.Range("H121").Value = 0.65
$ws.Range("J121").Value = 10.684931506849315
$ws.Range("K121").Value = 112.89863013698628
$ws.Range("L121").Value = 0.78466666666666596
$ws.Range("H122").Value = 0.65
$ws.Range("J122").Value = 16.027397260273972
$ws.Range("K122").Value = 118.34794520547946
$ws.Range("L122").Value = 0.74411111111111161
$ws.Range("H123").Value = 0.65
$ws.Range("J123").Value = 32.054794520547944
$ws.Range("K123").Value = 134.6958904109589
$ws.Range("L123").Value = 0.70355555555555516
$ws.Range("H124").Value = 0.65
$ws.Range("J124").Value = 48.082191780821923
$ws.Range("K124").Value = 151.04383561643834
$ws.Range("L124").Value = 0.69003703703703689
$ws.Range("H125").Value = 0.65
$ws.Range("J125").Value = 64.109589041095887
$ws.Range("K125").Value = 167.39178082191779
$ws.Range("L125").Value = 0.68327777777777765

